# regen sval data to filter save games
# Updates columns B (TB), C (d2S), D (K), E (IP), and G (sum) for rows 2-11
# on the active worksheet with the regenerated values. Column F (Win) and
# column A (dates) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.182878228561681,  1.65323645889881,   16.98373111632243,  0.4998867070740569, 22.31973251085698)
    3  = @(3.182878228561681,  1.65323645889881,   3.082599426703578,  6.48142807727062,   14.40014219143469)
    4  = @(3.182878228561681,  1.65323645889881,   0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    5  = @(3.182878228561681,  1.65323645889881,   0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    6  = @(3.182878228561681,  1.65323645889881,   3.082599426703578,  0.4998867070740569, 8.418600821238126)
    7  = @(3.182878228561681,  1.65323645889881,   0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    8  = @(3.182878228561681,  1.65323645889881,   0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    9  = @(0.1554434735375247, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569, 0.8605486643198037)
    10 = @(0.06328177979961902, 0.004309184025731883, 157.8057217802531, 6.48142807727062, 164.3547408213491)
    11 = @(0.7287194209349384, 0.3375848360084654, 16.98373111632243,  0.4998867070740569, 18.54992208033989)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
